$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.68650257129698
$ws.Range("D2").Value = 8.667714164980886
$ws.Range("E2").Value = 15.09451614514763
$ws.Range("F2").Value = 37.60032586645614
$ws.Range("G2").Value = 44.7401163386219
$ws.Range("H2").Value = 17.97516079517596
$ws.Range("J2").Value = 11.10524286567828
$ws.Range("L2").Value = 14.04779170350033
$ws.Range("B3").Value = 21.23538984763666
$ws.Range("D3").Value = 8.576648291965904
$ws.Range("E3").Value = 14.89523465863782
$ws.Range("F3").Value = 37.78666690255565
$ws.Range("G3").Value = 44.70362448172449
$ws.Range("H3").Value = 18.04869035474507
$ws.Range("J3").Value = 11.04065086285119
$ws.Range("L3").Value = 13.64899217655509
$ws.Range("B4").Value = 20.9558676951086
$ws.Range("D4").Value = 8.519833652220445
$ws.Range("E4").Value = 14.77187368869249
$ws.Range("F4").Value = 37.91666928339373
$ws.Range("G4").Value = 44.70560850360653
$ws.Range("H4").Value = 18.09974820118127
$ws.Range("J4").Value = 11.00156313243735
$ws.Range("L4").Value = 13.39863237776992
$ws.Range("B5").Value = 20.84146367270202
$ws.Range("D5").Value = 8.496466450573717
$ws.Range("E5").Value = 14.72139324218857
$ws.Range("F5").Value = 37.97353159999022
$ws.Range("G5").Value = 44.71251316598977
$ws.Range("H5").Value = 18.12203120849213
$ws.Range("J5").Value = 10.98578788225719
$ws.Range("L5").Value = 13.29535807194346
$ws.Range("B6").Value = 20.82244121038295
$ws.Range("D6").Value = 8.492573710785805
$ws.Range("E6").Value = 14.71299956391199
$ws.Range("F6").Value = 37.98320718319436
$ws.Range("G6").Value = 44.71402661264512
$ws.Range("H6").Value = 18.12582017511284
$ws.Range("J6").Value = 10.98317794726359
$ws.Range("L6").Value = 13.27813793047807
$ws.Range("B7").Value = 20.95432662245724
$ws.Range("D7").Value = 8.519519368272348
$ws.Range("E7").Value = 14.77119368558188
$ws.Range("F7").Value = 37.91742046352236
$ws.Range("G7").Value = 44.70567699517306
$ws.Range("H7").Value = 18.10004275236222
$ws.Range("J7").Value = 11.00134974874028
$ws.Range("L7").Value = 13.39724446807969
$ws.Range("B8").Value = 21.53156888854253
$ws.Range("D8").Value = 8.636508260332173
$ws.Range("E8").Value = 15.02603669517294
$ws.Range("F8").Value = 37.6613204080066
$ws.Range("G8").Value = 44.72245239986478
$ws.Range("H8").Value = 17.99928132082286
$ws.Range("J8").Value = 11.08285855523557
$ws.Range("L8").Value = 13.9115019775876
$ws.Range("B9").Value = 22.63745270749197
$ws.Range("D9").Value = 8.858288827774947
$ws.Range("E9").Value = 15.51603123333377
$ws.Range("F9").Value = 37.28429800972787
$ws.Range("G9").Value = 44.95005716379584
$ws.Range("H9").Value = 17.8490080276496
$ws.Range("J9").Value = 11.24681419500431
$ws.Range("L9").Value = 14.87107429028142
$ws.Range("B10").Value = 23.42649031783405
$ws.Range("D10").Value = 9.015931545344916
$ws.Range("E10").Value = 15.86772923196177
$ws.Range("F10").Value = 37.08561084824171
$ws.Range("G10").Value = 45.23690503117132
$ws.Range("H10").Value = 17.76801196724085
$ws.Range("J10").Value = 11.36921161687666
$ws.Range("L10").Value = 15.5398763351911
$ws.Range("B11").Value = 23.77888668752309
$ws.Range("D11").Value = 9.086364182415227
$ws.Range("E11").Value = 16.0254450090464
$ws.Range("F11").Value = 37.01263144562117
$ws.Range("G11").Value = 45.39335986630829
$ws.Range("H11").Value = 17.73766517537328
$ws.Range("J11").Value = 11.42519103154741
$ws.Range("L11").Value = 15.83511098279994
$ws.Range("B12").Value = 23.91127714585868
$ws.Range("D12").Value = 9.112840288403959
$ws.Range("E12").Value = 16.08480256062649
$ws.Range("F12").Value = 36.98753003667197
$ws.Range("G12").Value = 45.45632195269073
$ws.Range("H12").Value = 17.72711701192494
$ws.Range("J12").Value = 11.44642158411806
$ws.Range("L12").Value = 15.94552895789311
$ws.Range("B13").Value = 23.88281311767533
$ws.Range("D13").Value = 9.107147041514477
$ws.Range("E13").Value = 16.07203573380224
$ws.Range("F13").Value = 36.9928228730083
$ws.Range("G13").Value = 45.4425970610561
$ws.Range("H13").Value = 17.72934664564195
$ws.Range("J13").Value = 11.44184793757781
$ws.Range("L13").Value = 15.92181116450019
$ws.Range("B14").Value = 23.78980028677295
$ws.Range("D14").Value = 9.088546362628461
$ws.Range("E14").Value = 16.03033594409023
$ws.Range("F14").Value = 37.01051537067022
$ws.Range("G14").Value = 45.39846545845712
$ws.Range("H14").Value = 17.73677841213708
$ws.Range("J14").Value = 11.42693709521705
$ws.Range("L14").Value = 15.84422327345231
$ws.Range("B15").Value = 23.73268668034313
$ws.Range("D15").Value = 9.077127162214955
$ws.Range("E15").Value = 16.004744821552
$ws.Range("F15").Value = 37.02168349012705
$ws.Range("G15").Value = 45.37191681198547
$ws.Range("H15").Value = 17.7414537168392
$ws.Range("J15").Value = 11.41780766742662
$ws.Range("L15").Value = 15.79651620808397
$ws.Range("B16").Value = 23.40331886293708
$ws.Range("D16").Value = 9.011301891749142
$ws.Range("E16").Value = 15.85737301446342
$ws.Range("F16").Value = 37.09073317682564
$ws.Range("G16").Value = 45.22720155038188
$ws.Range("H16").Value = 17.77012675032693
$ws.Range("J16").Value = 11.36555834808991
$ws.Range("L16").Value = 15.52039328071836
$ws.Range("B17").Value = 23.19950370204453
$ws.Range("D17").Value = 8.970584382694174
$ws.Range("E17").Value = 15.76635464683401
$ws.Range("F17").Value = 37.1375735239909
$ws.Range("G17").Value = 45.1450645146399
$ws.Range("H17").Value = 17.78938772872356
$ws.Range("J17").Value = 11.33357429209502
$ws.Range("L17").Value = 15.34863009300971
$ws.Range("B18").Value = 23.0816652884636
$ws.Range("D18").Value = 8.947044934263475
$ws.Range("E18").Value = 15.71379155358881
$ws.Range("F18").Value = 37.16615123467374
$ws.Range("G18").Value = 45.10026725876005
$ws.Range("H18").Value = 17.80107730960734
$ws.Range("J18").Value = 11.31520681523676
$ws.Range("L18").Value = 15.24899324746869
$ws.Range("B19").Value = 23.0416661681664
$ws.Range("D19").Value = 8.939054648668911
$ws.Range("E19").Value = 15.69595944213468
$ws.Range("F19").Value = 37.17610722074195
$ws.Range("G19").Value = 45.08552012373673
$ws.Range("H19").Value = 17.80513989033121
$ws.Range("J19").Value = 11.30899319682431
$ws.Range("L19").Value = 15.21511586216892
$ws.Range("B20").Value = 23.22126408547167
$ws.Range("D20").Value = 8.974931313237132
$ws.Range("E20").Value = 15.7760659089889
$ws.Range("F20").Value = 37.13241769647516
$ws.Range("G20").Value = 45.15355510932926
$ws.Range("H20").Value = 17.78727404396779
$ws.Range("J20").Value = 11.33697613601047
$ws.Range("L20").Value = 15.36700249539161
$ws.Range("B21").Value = 23.81714989522145
$ws.Range("D21").Value = 9.094015218147575
$ws.Range("E21").Value = 16.04259442010239
$ws.Range("F21").Value = 37.00524963047065
$ws.Range("G21").Value = 45.41132732027985
$ws.Range("H21").Value = 17.73456984580611
$ws.Range("J21").Value = 11.43131597891505
$ws.Range("L21").Value = 15.86705083781354
$ws.Range("B22").Value = 24.20040248567215
$ws.Range("D22").Value = 9.170700677810578
$ws.Range("E22").Value = 16.21463587672051
$ws.Range("F22").Value = 36.93692341481859
$ws.Range("G22").Value = 45.60144090793676
$ws.Range("H22").Value = 17.70562799040234
$ws.Range("J22").Value = 11.49315614340208
$ws.Range("L22").Value = 16.18577819226134
$ws.Range("B23").Value = 23.99645567411907
$ws.Range("D23").Value = 9.129880373515281
$ws.Range("E23").Value = 16.12302339177628
$ws.Range("F23").Value = 36.97202748309454
$ws.Range("G23").Value = 45.49800173094974
$ws.Range("H23").Value = 17.7205683703766
$ws.Range("J23").Value = 11.46013765368379
$ws.Range("L23").Value = 16.01643312996597
$ws.Range("B24").Value = 23.21142827051629
$ws.Range("D24").Value = 8.972966471286721
$ws.Range("E24").Value = 15.7716761776735
$ws.Range("F24").Value = 37.13474351357461
$ws.Range("G24").Value = 45.14970895870773
$ws.Range("H24").Value = 17.78822772224484
$ws.Range("J24").Value = 11.33543809744562
$ws.Range("L24").Value = 15.35869909128308
$ws.Range("B25").Value = 22.34187976437507
$ws.Range("D25").Value = 8.799176042018562
$ws.Range("E25").Value = 15.38476587250175
$ws.Range("F25").Value = 37.37268040220974
$ws.Range("G25").Value = 44.86750373065513
$ws.Range("H25").Value = 17.88453583877779
$ws.Range("J25").Value = 11.20207752288216
$ws.Range("L25").Value = 14.61741554619385
